# Q3 Update - 2025
# UNHCR UN-LKA dataset refresh:
#  - shared "short-url" slug changes for every data row (gUe1ED -> y1mG0k)
#  - the last three rows of the 2024 country-of-origin breakdown are dropped
#  - rows 195-200 shift up to absorb the row below them (country-of-origin +
#    a few refreshed counts), and two further values are refreshed in place
#  - sheet shrinks from A1:V203 to A1:V200

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# 1) Refresh the short-url slug shared by every data row (row 2 .. last row)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value() = "y1mG0k"
}

# 2) Update the coo_id / coo_name / coo / coo_iso (F:I) and the refugee
#    count columns (N, O, Q, S) for rows 194-200 so each row absorbs the
#    data that used to live one row below it, plus a couple of refreshed
#    figures (N194, S200).
$ws.Cells.Item(194, 14).Value() = "27"

$ws.Cells.Item(195, 6).Value()  = "24"
$ws.Cells.Item(195, 7).Value()  = "Belarus"
$ws.Cells.Item(195, 8).Value()  = "BLR"
$ws.Cells.Item(195, 9).Value()  = "BLR"

$ws.Cells.Item(196, 6).Value()  = "69"
$ws.Cells.Item(196, 7).Value()  = "Palestinian"
$ws.Cells.Item(196, 8).Value()  = "GAZ"
$ws.Cells.Item(196, 9).Value()  = "PSE"
$ws.Cells.Item(196, 14).Value() = "5"
$ws.Cells.Item(196, 15).Value() = "0"

$ws.Cells.Item(197, 6).Value()  = "112"
$ws.Cells.Item(197, 7).Value()  = "Sri Lanka"
$ws.Cells.Item(197, 8).Value()  = "LKA"
$ws.Cells.Item(197, 9).Value()  = "LKA"
$ws.Cells.Item(197, 14).Value() = "0"
$ws.Cells.Item(197, 17).Value() = "4806"

$ws.Cells.Item(198, 6).Value()  = "135"
$ws.Cells.Item(198, 7).Value()  = "Myanmar"
$ws.Cells.Item(198, 8).Value()  = "MYA"
$ws.Cells.Item(198, 9).Value()  = "MMR"
$ws.Cells.Item(198, 14).Value() = "112"
$ws.Cells.Item(198, 15).Value() = "117"
$ws.Cells.Item(198, 17).Value() = "0"

$ws.Cells.Item(199, 6).Value()  = "147"
$ws.Cells.Item(199, 7).Value()  = "Pakistan"
$ws.Cells.Item(199, 8).Value()  = "PAK"
$ws.Cells.Item(199, 9).Value()  = "PAK"
$ws.Cells.Item(199, 14).Value() = "88"
$ws.Cells.Item(199, 15).Value() = "146"

$ws.Cells.Item(200, 6).Value()  = "216"
$ws.Cells.Item(200, 7).Value()  = "Stateless"
$ws.Cells.Item(200, 8).Value()  = "STA"
$ws.Cells.Item(200, 9).Value()  = "XXA"
$ws.Cells.Item(200, 14).Value() = "0"
$ws.Cells.Item(200, 15).Value() = "0"
$ws.Cells.Item(200, 19).Value() = "229"

# 3) Drop the trailing three rows (201-203) that the shift above makes
#    redundant - delete bottom-up so row numbers stay stable.
$ws.Rows.Item(203).Delete()
$ws.Rows.Item(202).Delete()
$ws.Rows.Item(201).Delete()
